# The two rows (2 and 3) describe the same sighting/location pair but
# were re-ordered: row 2's Id/Ost/Nord need to become row 3's original
# values and vice versa (everything else in the rows is identical).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = $ws.Range("A2").Value()
$q2 = $ws.Range("Q2").Value()
$r2 = $ws.Range("R2").Value()

$a3 = $ws.Range("A3").Value()
$q3 = $ws.Range("Q3").Value()
$r3 = $ws.Range("R3").Value()

$ws.Range("A2").Value = $a3
$ws.Range("Q2").Value = $q3
$ws.Range("R2").Value = $r3

$ws.Range("A3").Value = $a2
$ws.Range("Q3").Value = $q2
$ws.Range("R3").Value = $r2
